$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# A new September entry was logged (timestamped 2024-09-09 11:21:43).
# The tracker sheet keeps entries sorted with newest first, so the new
# entry is inserted as a fresh row above the prior top-of-September row
# (row 35). That shifts everything below (the rest of September's R/S
# entries, and the trailing August "hdfc" P/Q entries, plus the
# "Broadband" label) down by exactly one row.
$ws.Rows.Item(35).Insert()

$ws.Cells.Item(35, 18).Value = "corporate internet share"
$ws.Cells.Item(35, 19).Value = "2024-09-09 11:21:43"
